# Append the 2025-03-13 price row to every sheet in the workbook.
# Each worksheet tracks a "Date"/"Price" series in columns A/B starting at
# row 2; the most recent row (row 11, 2025-03-12) is duplicated forward to
# a new row 12 dated 2025-03-13, carrying the same price as the previous
# day (matches the source diff, which shows every sheet's new row repeating
# the prior day's value except USD_CNY, whose own latest value is reused).
#
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the workbook's existing convention of keeping the
# Date/Price columns as text) rather than auto-converting the date-shaped
# strings into date serials or the comma-grouped numbers into formatted
# numerics.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-13"

$prices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,289"
    "Silver Busbar front-side"   = "7,917"
    "Silver finger front-side"   = "7,967"
    "USD_CNY"                    = "7.2506"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    $price = $prices[$name]

    $ws.Range("A12").Value = "'" + $newDate
    $ws.Range("B12").Value = "'" + $price
}
